$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows to match repulled data
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -5
$ws.Range("F8").Value = 6
$ws.Range("F9").Value = -1
$ws.Range("F12").Value = 6
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = -6
$ws.Range("F22").Value = 9
$ws.Range("F23").Value = -2
